$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "05/05/2021 01:39:46"
$ws.Cells.Item(2, 2).Value = 15.666
$ws.Cells.Item(3, 1).Value = "05/05/2021 01:40:48"
$ws.Cells.Item(3, 2).Value = 14.31
$ws.Cells.Item(4, 1).Value = "05/05/2021 01:41:50"
$ws.Cells.Item(4, 2).Value = 15.114
$ws.Cells.Item(5, 1).Value = "05/05/2021 01:42:55"
$ws.Cells.Item(5, 2).Value = 15.114
$ws.Cells.Item(6, 1).Value = "05/05/2021 01:44:00"
$ws.Cells.Item(6, 2).Value = 14.95
$ws.Cells.Item(7, 1).Value = "05/05/2021 01:45:04"
$ws.Cells.Item(7, 2).Value = 14.798
$ws.Cells.Item(8, 1).Value = "05/05/2021 01:46:09"
$ws.Cells.Item(8, 2).Value = 18.395
$ws.Cells.Item(9, 1).Value = "05/05/2021 01:51:31"
$ws.Cells.Item(9, 2).Value = 38.542
$ws.Cells.Item(10, 1).Value = "05/05/2021 01:52:37"
$ws.Cells.Item(10, 2).Value = 19.865
$ws.Cells.Item(11, 1).Value = "05/05/2021 01:53:43"
$ws.Cells.Item(11, 2).Value = 37.938
$ws.Cells.Item(12, 1).Value = "05/05/2021 01:54:50"
$ws.Cells.Item(12, 2).Value = 19.284
$ws.Cells.Item(13, 1).Value = "05/05/2021 02:00:08"
$ws.Cells.Item(13, 2).Value = 18.601
$ws.Cells.Item(14, 1).Value = "05/05/2021 02:01:17"
$ws.Cells.Item(14, 2).Value = 18.398
$ws.Cells.Item(15, 1).Value = "05/05/2021 02:06:33"
$ws.Cells.Item(15, 2).Value = 105687.46
$ws.Cells.Item(16, 1).Value = "05/05/2021 02:07:38"
$ws.Cells.Item(16, 2).Value = 28.264
$ws.Cells.Item(17, 1).Value = "05/05/2021 02:12:49"
$ws.Cells.Item(17, 2).Value = 14.396
$ws.Cells.Item(18, 1).Value = "05/05/2021 02:13:53"
$ws.Cells.Item(18, 2).Value = 14.638
$ws.Cells.Item(19, 1).Value = "05/05/2021 02:14:57"
$ws.Cells.Item(19, 2).Value = 16.392
$ws.Cells.Item(20, 1).Value = "05/05/2021 02:16:00"
$ws.Cells.Item(20, 2).Value = 14.32
$ws.Cells.Item(21, 1).Value = "05/05/2021 02:17:04"
$ws.Cells.Item(21, 2).Value = 14.09
$ws.Cells.Item(22, 1).Value = "05/05/2021 02:18:07"
$ws.Cells.Item(22, 2).Value = 16.332
$ws.Cells.Item(23, 1).Value = "05/05/2021 02:19:12"
$ws.Cells.Item(23, 2).Value = 15.668
$ws.Cells.Item(24, 1).Value = "05/05/2021 02:20:16"
$ws.Cells.Item(24, 2).Value = 12.683
$ws.Cells.Item(25, 1).Value = "05/05/2021 02:21:21"
$ws.Cells.Item(25, 2).Value = 19.073
$ws.Cells.Item(26, 1).Value = "05/05/2021 02:21:31"
$ws.Cells.Item(26, 2).Value = 19.073
$ws.Cells.Item(27, 1).Value = "05/05/2021 02:21:40"
$ws.Cells.Item(27, 2).Value = 19.073
$ws.Cells.Item(28, 1).Value = "05/05/2021 02:21:48"
$ws.Cells.Item(28, 2).Value = 19.073
$ws.Cells.Item(29, 1).Value = "05/05/2021 02:21:58"
$ws.Cells.Item(29, 2).Value = 20.152
$ws.Cells.Item(30, 1).Value = "05/05/2021 02:22:08"
$ws.Cells.Item(30, 2).Value = 20.152
$ws.Cells.Item(31, 1).Value = "05/05/2021 02:22:18"
$ws.Cells.Item(31, 2).Value = 20.152
$ws.Cells.Item(32, 1).Value = "05/05/2021 02:22:30"
$ws.Cells.Item(32, 2).Value = 20.152
$ws.Cells.Item(33, 1).Value = "05/05/2021 02:22:40"
$ws.Cells.Item(33, 2).Value = 20.152
$ws.Cells.Item(34, 1).Value = "05/05/2021 02:22:49"
$ws.Cells.Item(34, 2).Value = 20.152

# Delete rows 35 through 53 which no longer have data (dimension shrinks from A1:B53 to A1:B34)
$ws.Range("A35:B53").EntireRow.Delete() | Out-Null
